$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Append the new "Vendredi 05 juin 2020" journal entry block at the
#    very end of the document body (before the final sectPr), exactly
#    mirroring the target OOXML structure via a single InsertXML call.
# ------------------------------------------------------------------
$contentRange = $d.Content
$contentRange.Collapse(0)

$newBlockXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:br w:type="page"/>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading2"/>
            </w:pPr>
            <w:r>
              <w:lastRenderedPageBreak/>
              <w:t>Vendredi 05 juin 2020 – 9</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:vertAlign w:val="superscript"/>
              </w:rPr>
              <w:t>ème</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> jour</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading3"/>
            </w:pPr>
            <w:r>
              <w:t>08h00</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>J’ai continué à peaufiner mon site pendant la matinée. J’ai eu quelques soucis en JQuery pour afficher les images dynamiquement après le choix de l’utilisateur. J’avais réussi assez facilement à afficher une image, mais lorsqu’il s’agissait de plusieurs médias, ça ne marchait pas. Cette fonctionnalité m’a donc pris un certain temps à coder.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading3"/>
            </w:pPr>
            <w:r>
              <w:t>12h00</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>J’ai documenté tous les changements faits sur mon site. J’ai également vérifié qu’il n’y avait pas d’erreur dans ma documentation technique.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading3"/>
            </w:pPr>
            <w:r>
              <w:t>13h00</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>J’ai commencé le manuel utilisateur. Le manuel utilisateur se trouve sur le site Web. C’est une page d’aide que les utilisateurs pourront consulter en cas de problèmes.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$contentRange.InsertXML($newBlockXml)

# ------------------------------------------------------------------
# 2) Update the cached field result in the header from "4 juin 2020"
#    to "6 juin 2020" (the TIME field itself is left untouched).
# ------------------------------------------------------------------
$headerRange = $d.Sections.Item(1).Headers.Item(1).Range
$headerRange.Find.Execute("4 juin 2020", $true, $false, $false, $false, $false, `
                           $true, 1, $false, "6 juin 2020", 2)
